$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("adminMentor", "abc123", "Invalid"),
    @("studentAdmin", "xyz123", "Invalid"),
    @("sangeeta", "sangeeta123", "Valid"),
    @("sweetapal", "sweeta123", "Valid")
)

$r = 4
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $r++
}

$ws.Range("A2:C2").Copy()
$ws.Range("A4:C7").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("C6").Select()
